# Update betting-odds values on Sheet1 for rows 4, 5 and 7
# (weekly refresh of FlashScore odds data).

$wb = $excel.ActiveWorkbook
$ws = $wb.Sheets.Item("Sheet1")

# ---- Row 4 ----
$ws.Range("G4").Value  = 2.72
$ws.Range("I4").Value  = 2.65
$ws.Range("J4").Value  = 3.35
$ws.Range("N4").Value  = 6.1
$ws.Range("O4").Value  = 1.39
$ws.Range("P4").Value  = 2.75
$ws.Range("Q4").Value  = 2.15
$ws.Range("R4").Value  = 1.65
$ws.Range("U4").Value  = 1.8
$ws.Range("V4").Value  = 1.9
$ws.Range("W4").Value  = 7.8
$ws.Range("AC4").Value = 6.1
$ws.Range("AD4").Value = 5.5
$ws.Range("AG4").Value = 7.5
$ws.Range("AH4").Value = 13
$ws.Range("AI4").Value = 9.75
$ws.Range("AK4").Value = 24
$ws.Range("AL4").Value = 35
$ws.Range("AO4").Value = 15.5
$ws.Range("AS4").Value = 300
$ws.Range("AW4").Value = 4.6
$ws.Range("AY4").Value = 22
$ws.Range("BA4").Value = 100

# ---- Row 5 ----
$ws.Range("G5").Value  = 8
$ws.Range("H5").Value  = 4.15
$ws.Range("I5").Value  = 1.36
$ws.Range("J5").Value  = 7.7
$ws.Range("K5").Value  = 2.27
$ws.Range("L5").Value  = 1.85
$ws.Range("M5").Value  = 1.06
$ws.Range("N5").Value  = 7.4
$ws.Range("O5").Value  = 1.28
$ws.Range("P5").Value  = 3.35
$ws.Range("Q5").Value  = 1.85
$ws.Range("R5").Value  = 1.9
$ws.Range("S5").Value  = 1.39
$ws.Range("T5").Value  = 2.77
$ws.Range("W5").Value  = 18
$ws.Range("X5").Value  = 55
$ws.Range("Y5").Value  = 27
$ws.Range("Z5").Value  = 250
$ws.Range("AC5").Value = 7.4
$ws.Range("AD5").Value = 8.5
$ws.Range("AE5").Value = 23
$ws.Range("AH5").Value = 5.9
$ws.Range("AJ5").Value = 8.5
$ws.Range("AK5").Value = 12
$ws.Range("AN5").Value = 9.25
$ws.Range("AR5").Value = 500
$ws.Range("AT5").Value = 2.77
$ws.Range("AX5").Value = 6.2
$ws.Range("AZ5").Value = 17.5

# ---- Row 7 ----
$ws.Range("H7").Value  = 3.65
$ws.Range("I7").Value  = 1.57
$ws.Range("K7").Value  = 2.15
$ws.Range("L7").Value  = 2.15
$ws.Range("M7").Value  = 1.06
$ws.Range("N7").Value  = 7.1
$ws.Range("Q7").Value  = 1.95
$ws.Range("R7").Value  = 1.8
$ws.Range("S7").Value  = 1.42
$ws.Range("T7").Value  = 2.67
$ws.Range("Z7").Value  = 100
$ws.Range("AC7").Value = 7.1
$ws.Range("AD7").Value = 7.2
$ws.Range("AG7").Value = 6.2
$ws.Range("AJ7").Value = 11.25
$ws.Range("AN7").Value = 6.8
$ws.Range("AT7").Value = 2.67
$ws.Range("AU7").Value = 8
$ws.Range("AV7").Value = 80
$ws.Range("AW7").Value = 3.35
$ws.Range("AX7").Value = 7.8
$ws.Range("AY7").Value = 19
$ws.Range("AZ7").Value = 26
$ws.Range("BA7").Value = 60
$ws.Range("BB7").Value = 300
